$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Translations"
$ws.Name = "Translations"

$guid = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"

# Row 1 - header: a new "Entity Id" column is introduced at column A, pushing
# the old Type/Index/Id headers one column to the right and dropping "Id".
$ws.Range("A1").Value = "Entity Id"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Index"
$ws.Range("D1").Value = "Original"
$ws.Range("E1").Value = "Translation"

# Row 2 - Title entity (no Index value)
$ws.Range("A2").Value = $guid
$ws.Range("B2").Value = "Title"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "Orig"
$ws.Range("E2").Value = "title"

# Row 3 - ValidationMessage entity (Index = 1)
$ws.Range("A3").Value = $guid
$ws.Range("B3").Value = "ValidationMessage"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Orig"
$ws.Range("E3").Value = "validation message"

# Row 4 - Instruction entity (no Index value)
$ws.Range("A4").Value = $guid
$ws.Range("B4").Value = "Instruction"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "Orig"
$ws.Range("E4").Value = "instruction"

# Row 5 - OptionTitle entity (Index = 2)
$ws.Range("A5").Value = $guid
$ws.Range("B5").Value = "OptionTitle"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "Orig"
$ws.Range("E5").Value = "option"

# Column widths: A/B/C now hold the new data and are auto-fit; column C
# shrinks a lot (now holding short Index numbers rather than the long Id
# guid). Column E's width is left untouched (it already has the right custom
# width from before the edit).
$ws.Columns.Item(1).ColumnWidth = 42.5
$ws.Columns.Item(2).ColumnWidth = 17.33
$ws.Columns.Item(3).ColumnWidth = 5.2

# Move the active selection to E10, matching the saved workbook state.
$ws.Range("E10").Select() | Out-Null
